# "Generate Report for Archive"
#
# The localization status report is regenerated: the "Status" value that was
# "Ready for handoff" is now "In Translation" (it shows up on the Overview
# sheet's per-language summary columns as well as on each language sheet's
# "Status" column). Because the status text got shorter, the columns that
# display it were re-sized to fit the new content.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status text ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-language sheets: column C ("Status") holds the status text ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
